# Update the R6_Legacy_compare workbook to use the backslash-escaped
# "$" sigils in the R6 method-call examples on the RLcomp_valid sheet
# (consistent with the other sheets), and restore the single-cell
# selection on that sheet.

$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("RLcomp_valid")

$ws.Range("B6").Value = "PM_result\`$step()"
$ws.Range("B5").Value = "PM_valid\`$plot()"
$ws.Range("B3").Value = "PM_result\`$op\`$plot(resid = T,…)"

# Shrink the selection on RLcomp_valid from B4:C4 down to just B4, without
# permanently changing which sheet tab is active in the workbook.
$ws.Range("B4").Select() | Out-Null
$originalActive.Activate() | Out-Null
